# results of GoodSplit with TeacherRL
# Append a new results row (row 36) to the DFA_R97_sched4 sheet describing
# a GoodSplit run (TeacherRL, maxLen:3) and select the new row's cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A36").Value = 1
$ws.Range("B36").Value = 350368
$ws.Range("C36").Value = 350368
$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 4772106
$ws.Range("F36").Value = "../data/experiments/DFA_R97_sched4.fsm"

# Note: write order matters for shared-string table append order.
$ws.Range("L36").Value = "Correct: 1, reset: 350368,      OQ: 350368,     EQ: 0,  symbols: 4772106,"
$ws.Range("G36").Value = "GoodSplit"
$ws.Range("H36").Value = "maxLen:3"
$ws.Range("I36").Value = "TeacherRL"

[void]$ws.Range("H36").Select()
